$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "27.040.66"
$ws.Range("D3").Value = "1.729.31"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.42"
$ws.Range("E5").Value = "  -5.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4863"
$ws.Range("E7").Value = "  +4.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3490"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.40"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.99"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.889"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "1.726.70"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.854"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.11"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06395"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.726"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "27.092.39"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.96"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.074"
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.33"
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "1.925.37"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.075"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.99"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.045"
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09353"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.644"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.401"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05926"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02185"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.99"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.427"
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1999"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.762"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9997"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5994"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.118"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.541"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.85"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.580"
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5630"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.12"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.852"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.105"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06653"
$ws.Range("E51").Value = "  -2.17%  "
